$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data to append after the last existing row (row 208)
$newRows = @(
    @{ Title = "病娇男孩的精分日记"; Link = "https://pan.quark.cn/s/992407d06077" },
    @{ Title = "桔梗物语"; Link = "https://pan.quark.cn/s/be629be38438" },
    @{ Title = "桐花中路私立协济医院怪谈"; Link = "https://pan.quark.cn/s/3acf12cb40ce" }
)

$startRow = 209
$category = $ws.Cells.Item(208, 3).Value

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $r - 2 + 1
    $ws.Cells.Item($r, 2).Value = $row.Title
    $ws.Cells.Item($r, 3).Value = $category
    $ws.Cells.Item($r, 4).Value = $row.Link

    $ws.Rows.Item($r).RowHeight = 15.5
}

$ws.Range("C213").Select()
